# I0 and IF added
# Adds two new columns (I: "I0", J: "IF") to the existing data table on the
# active sheet, mirroring the header style used by the other header cells
# (copied from H1) and filling in the per-row numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the formatting of the existing header cells (bold, bordered,
# centered) by copying H1's style onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data ----------------------------------------------------------------
$colI = @(1,1,1,1,1,7,6,9,2,1,1,2,1,2,8,3,3,4,3,7,6,10,8,7,7,7,9,6,4,4,1,1,1,1,1,1,1,1,1,1,1,1,1)
$colJ = @(7,7,7,6,5,9,7,9,6,5,5,5,8,6,8,8,7,7,6,7,8,10,8,7,7,7,9,8,6,6,6,4,6,5,6,5,5,5,6,5,4,3,2)

for ($r = 2; $r -le 44; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $colI[$idx]
    $ws.Cells.Item($r, 10).Value = $colJ[$idx]
}
